# Add "Partigrupper" worksheet: a PartiGruppeKort -> PartiGruppe mapping table,
# matching the shape/position of the existing "Lande"/"Ophold" mapping sheets.
# Cell values are written in the exact order the source data was authored so the
# resulting sharedStrings.xml table comes out in the same sequence.

$wb = $excel.ActiveWorkbook

$writes = @(
    @{Cell='A2'; Value='S'; Styled=$false},
    @{Cell='A3'; Value='DF'; Styled=$false},
    @{Cell='A4'; Value='V'; Styled=$false},
    @{Cell='A5'; Value='RV'; Styled=$false},
    @{Cell='A6'; Value='SF'; Styled=$false},
    @{Cell='A7'; Value='KF'; Styled=$false},
    @{Cell='A8'; Value='EL'; Styled=$false},
    @{Cell='A9'; Value='M'; Styled=$false},
    @{Cell='A10'; Value='DD'; Styled=$false},
    @{Cell='A11'; Value='LA'; Styled=$false},
    @{Cell='A12'; Value='Udlændinge- Og Integrationsministeren'; Styled=$false},
    @{Cell='A13'; Value='ALT'; Styled=$false},
    @{Cell='A14'; Value='NB'; Styled=$false},
    @{Cell='A15'; Value='Integrationsministeren'; Styled=$false},
    @{Cell='A16'; Value='Justitsministeren'; Styled=$false},
    @{Cell='A17'; Value='Den Fg. Formand'; Styled=$false},
    @{Cell='A18'; Value='UFG'; Styled=$false},
    @{Cell='A19'; Value='Udlændinge-, Integrations- Og Boligministeren'; Styled=$false},
    @{Cell='A20'; Value='Fg. Formand'; Styled=$false},
    @{Cell='A21'; Value='KD'; Styled=$false},
    @{Cell='A22'; Value='FG'; Styled=$false},
    @{Cell='A23'; Value='Fjerde Næstformand'; Styled=$false},
    @{Cell='A24'; Value='Første Næstformand'; Styled=$false},
    @{Cell='A25'; Value='Tredje Næstformand'; Styled=$false},
    @{Cell='A26'; Value='Anden Næstformand'; Styled=$false},
    @{Cell='A1'; Value='PartiGruppeKort'; Styled=$false},
    @{Cell='B1'; Value='PartiGruppe'; Styled=$false},
    @{Cell='B12'; Value='Udlændingeministeren'; Styled=$false},
    @{Cell='B17'; Value='Folketingets formand'; Styled=$false},
    @{Cell='B2'; Value='Socialdemokratiet (S)'; Styled=$true},
    @{Cell='B4'; Value='Venstre (V)'; Styled=$true},
    @{Cell='B3'; Value='Dansk Folkeparti (DF)'; Styled=$true},
    @{Cell='B6'; Value='Socialistisk Folkeparti (SF)'; Styled=$true},
    @{Cell='B5'; Value='Det Radikale Venstre (RV)'; Styled=$true},
    @{Cell='B7'; Value='Det Konservative Folkeparti (KF)'; Styled=$true},
    @{Cell='B8'; Value='Enhedslisten (EL)'; Styled=$true},
    @{Cell='B11'; Value='Liberal Alliance (LA)'; Styled=$true},
    @{Cell='B21'; Value='Kristendemokraterne (KD)'; Styled=$true},
    @{Cell='B18'; Value='Uden for folketingsgrupperne (UFG)'; Styled=$true},
    @{Cell='B13'; Value='Alternativet (ALT)'; Styled=$true},
    @{Cell='B14'; Value='Nye Borgerlige (NB)'; Styled=$true},
    @{Cell='B22'; Value='Frie Grønne, Danmarks Nye Venstrefløjsparti (FG)'; Styled=$true},
    @{Cell='B9'; Value='Moderaterne (M)'; Styled=$true},
    @{Cell='B10'; Value='Danmarksdemokraterne - Inger Støjberg (DD)'; Styled=$true},
    @{Cell='B15'; Value='Udlændingeministeren'; Styled=$false},
    @{Cell='B16'; Value='Udlændingeministeren'; Styled=$false},
    @{Cell='B19'; Value='Udlændingeministeren'; Styled=$false},
    @{Cell='B20'; Value='Folketingets formand'; Styled=$false},
    @{Cell='B23'; Value='Folketingets formand'; Styled=$false},
    @{Cell='B24'; Value='Folketingets formand'; Styled=$false},
    @{Cell='B25'; Value='Folketingets formand'; Styled=$false},
    @{Cell='B26'; Value='Folketingets formand'; Styled=$false}
)

# New sheet goes after the last existing sheet ("Ophold") and becomes the active/selected tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Partigrupper"

foreach ($w in $writes) {
    $cell = $ws.Range($w.Cell)
    $cell.Value = $w.Value
    if ($w.Styled) {
        $cell.Font.Name = "Calibri"
        $cell.Font.Color = 0
    }
}

$ws.Columns.Item(1).ColumnWidth = 38.33203125
$ws.Columns.Item(2).ColumnWidth = 41.44140625

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:B26"), $null, 1)
$lo.Name = "Partigrupper"
$lo.TableStyle = "TableStyleMedium4"

$ws.PageSetup.Orientation = 1

$ws.Range("C5").Select()
